$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.885.03'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.639.52'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = '''217.14'
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  +0.99%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  +3.33%  '
$ws.Range('D11').Value = '''0.0845'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '1.868.46'
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = '1.644.74'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D16').Value = '''67.04'
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').Value = '26.868.65'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = '''218.27'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '''6.72'
$ws.Range('E21').Value = '  +2.08%  '
$ws.Range('D22').Value = '''4.40'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = '''2.43'
$ws.Range('E23').Value = '  +2.87%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '''147.11'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').Value = '''7.26'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = '1.264.41'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('D38').Value = '''0.839'
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').Value = '''0.810'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('D43').Value = '1.779.44'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '''62.17'
$ws.Range('E44').Value = '  +1.59%  '
$ws.Range('D46').Value = '''91.94'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('D47').Value = '''1.60'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0512'
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''7.67'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.0961'
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '''1.00'
$ws.Range('E51').Value = '  -0.43%  '
